$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D30 gets the same "Description" text as D4 ("Fixed panel, \nHighlight active navigation item.")
$ws.Range("D30").Value = $ws.Range("D4").Value2

# D31 gets the same "Description" text as D5 ("implement houses page with house card list...")
$ws.Range("D31").Value = $ws.Range("D5").Value2

# Adjust row heights to fit the new/adjusted wrapped text
$ws.Rows.Item(5).RowHeight = 46
$ws.Rows.Item(30).RowHeight = 34
$ws.Rows.Item(31).RowHeight = 45

# Extend the SUMIF total range down through row 31 (previously E4:E29)
$ws.Range("B34").Formula = '=SUMIF(E4:E31,"<>x",B4:B31)'

# Update the view's scroll position and selection
$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 1
$ws.Range("D39").Select() | Out-Null
